$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# NOTE: whenever a new cell value looks like a number or a date (e.g.
# "20", "5", "12/12/2023") we pre-format the cell as Text ("@") so Excel
# stores it verbatim as a string instead of silently converting it to a
# number/date value.
# ----------------------------------------------------------------------

# ==========================================================================
# Sheet "Alunos": update H2, append rows 8 and 9
# ==========================================================================
$wsAlunos = $wb.Worksheets.Item("Alunos")

$wsAlunos.Range("H2").Value = "[7.14, 5.69, 5.69]"

# Row 8
$wsAlunos.Range("A8:H9").NumberFormat = "@"

$wsAlunos.Range("A8").Value = "ra054"
$wsAlunos.Range("B8").Value = "ra054"
$wsAlunos.Range("C8").Value = "daniel"
$wsAlunos.Range("D8").Value = "20"
$wsAlunos.Range("E8").Value = "da@da.com"
$wsAlunos.Range("F8").Value = "[]"
$wsAlunos.Range("G8").Value = "[]"
$wsAlunos.Range("H8").Value = ""

# Row 9
$wsAlunos.Range("A9").Value = "ra432"
$wsAlunos.Range("B9").Value = "ra432"
$wsAlunos.Range("C9").Value = "junn"
$wsAlunos.Range("D9").Value = "3"
$wsAlunos.Range("E9").Value = "dani@dan.com"
$wsAlunos.Range("F9").Value = "[]"
$wsAlunos.Range("G9").Value = "[]"
$wsAlunos.Range("H9").Value = ""

# ==========================================================================
# Sheet "Turmas": update E2 (append new ciclo entry), append rows 8 and 9
# ==========================================================================
$wsTurmas = $wb.Worksheets.Item("Turmas")

$wsTurmas.Range("E2").Value = '[{''id'': ''1'', ''nome'': ''C1'', ''data_de_inicio'': ''12'', ''data_de_fim'': ''21'', ''peso_da_nota'': ''3''}, {''id'': ''2'', ''nome'': ''C2'', ''data_de_inicio'': ''23'', ''data_de_fim'': ''32'', ''peso_da_nota'': ''4''}, {''id'': ''c9878959'', ''nome'': ''C3'', ''data_de_inicio'': ''12/12/2024'', ''data_de_fim'': ''13/12/2024'', ''peso_da_nota'': ''6''}, {''id'': ''c4088151'', ''nome'': ''aba'', ''data_de_inicio'': ''02/10/2024'', ''data_de_fim'': ''01/11/2024'', ''peso_da_nota'': ''5''}, {''id'': ''c4238120'', ''nome'': ''c4'', ''data_de_inicio'': ''27/11/2023'', ''data_de_fim'': ''28/11/2023'', ''peso_da_nota'': ''6''}, {''id'': ''c1773595'', ''nome'': ''aaaaa'', ''data_de_inicio'': ''27/11/2023'', ''data_de_fim'': ''27/12/2023'', ''peso_da_nota'': ''7''}, {''id'': ''c5717717'', ''nome'': ''3'', ''data_de_inicio'': ''12/12/2023'', ''data_de_fim'': ''12/12/2023'', ''peso_da_nota'': ''6''}, {''id'': ''c6884323'', ''nome'': ''c9878959'', ''data_de_inicio'': ''28/11/2023'', ''data_de_fim'': ''29/11/2023'', ''peso_da_nota'': ''6''}, {''id'': ''c5390957'', ''nome'': ''aaaaaaaa'', ''data_de_inicio'': ''29/11/2023'', ''data_de_fim'': ''30/11/2023'', ''peso_da_nota'': ''7''}, {''id'': ''C882'', ''nome'': ''nadalete '', ''data_de_inicio'': ''27/11/2023'', ''data_de_fim'': ''27/12/2023'', ''peso_da_nota'': ''5''}]'

# Rows 8 and 9
$wsTurmas.Range("A8:E9").NumberFormat = "@"

$wsTurmas.Range("A8").Value = "t3093321"
$wsTurmas.Range("B8").Value = "t3093321"
$wsTurmas.Range("C8").Value = "app"
$wsTurmas.Range("D8").Value = "20/12/2023"
$wsTurmas.Range("E8").Value = ""

$wsTurmas.Range("A9").Value = "t5300480"
$wsTurmas.Range("B9").Value = "t5300480"
$wsTurmas.Range("C9").Value = "133"
$wsTurmas.Range("D9").Value = "12/12/2023"
$wsTurmas.Range("E9").Value = '[{''id'': ''C039'', ''nome'': ''java'', ''data_de_inicio'': ''12/12/2024'', ''data_de_fim'': ''13/12/2024'', ''peso_da_nota'': ''6''}]'

# ==========================================================================
# Sheet "Ciclos": insert a new row 11 (shifting the old row 11 down to 12)
# and append a new row 13
# ==========================================================================
$wsCiclos = $wb.Worksheets.Item("Ciclos")

$wsCiclos.Rows.Item(11).Insert()

$wsCiclos.Range("A11:F11").NumberFormat = "@"
$wsCiclos.Range("A11").Value = "C882"
$wsCiclos.Range("B11").Value = "C882"
$wsCiclos.Range("C11").Value = "nadalete "
$wsCiclos.Range("D11").Value = "27/11/2023"
$wsCiclos.Range("E11").Value = "27/12/2023"
$wsCiclos.Range("F11").Value = "5"

# New row 13 (after old row 11 shifted down to row 12)
$wsCiclos.Range("A13:F13").NumberFormat = "@"
$wsCiclos.Range("A13").Value = "C039"
$wsCiclos.Range("B13").Value = "C039"
$wsCiclos.Range("C13").Value = "java"
$wsCiclos.Range("D13").Value = "12/12/2024"
$wsCiclos.Range("E13").Value = "13/12/2024"
$wsCiclos.Range("F13").Value = "6"

# ==========================================================================
# Sheet "Grupos": append row 7
# ==========================================================================
$wsGrupos = $wb.Worksheets.Item("Grupos")

$wsGrupos.Range("A7").Value = "G435"
$wsGrupos.Range("B7").Value = "G435"
$wsGrupos.Range("C7").Value = "c"
$wsGrupos.Range("D7").Value = "['ra054']"
